# Add a new row (row 3) to the flowcytometry "data_info" sheet that carries
# the French description / enum metadata for each column: a short code for
# the operator, sample id, date, operating mode, critical apparatus,
# critical product and raw-data storage location columns (A:G). The
# remaining columns (H:M) have no description text for this row, so they
# are touched (without giving them real content) purely so the row keeps
# the same shape/width as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"

# H3:M3 stay blank (no description) but are still part of row 3 — touch
# them (harmless no-op formatting) so the cells exist on the row.
$ws.Range("H3:M3").Font.Bold = $false
